# Gantt chart update: insert a new date column before column H (shifting
# the existing H:J columns to I:K), fill in the new milestone date, and
# leave the selection where the author left it when they saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at H - this shifts the old H,I,J columns (and their
# formatting) one column to the right, matching the new K3=45999 etc.
$ws.Columns("H:H").Insert() | Out-Null

# New milestone date for the inserted column (2025-11-03)
$ws.Range("H3").Value = 45964

# The author's last selection before saving.
$ws.Range("H21").Select() | Out-Null
